$wb = $excel.ActiveWorkbook

# --- Keywords sheet: add the two new keyword rows (9 and 10) ---
$wsKeywords = $wb.Worksheets.Item("Keywords")
$wsKeywords.Range("A9").Value = "120: KARTE / MAP"
$wsKeywords.Range("A10").Value = "121: KREATUREN DER SEE / SEA CREATURES"

# --- GlobalVars sheet: add the new global var row (21) ---
$wsGlobalVars = $wb.Worksheets.Item("GlobalVars")
$wsGlobalVars.Range("A21").Value = "238: You got the journal (log) from Torle"

# --- Items sheet: add the new item row (7) - Torle's journal ---
$wsItems = $wb.Worksheets.Item("Items")
$wsItems.Range("A7").Value = 408
$wsItems.Range("B7").Value = "Torle's Logbuch / Torle's Journal"
$wsItems.Range("D7").Value = "Contains the coordinates and or hints to all sea creatures"

# --- ObjectTexts sheet: add the new object text row (3) ---
$wsObjectTexts = $wb.Worksheets.Item("ObjectTexts")
$wsObjectTexts.Range("A3").Value = 2
$wsObjectTexts.Range("B3").Value = 3
$wsObjectTexts.Range("C3").Value = "Text of Torle's journal"

# Items C7 reuses the existing "Text Scroll" shared string, set last so it
# does not introduce a new shared-string entry out of order.
$wsItems.Range("C7").Value = "Text Scroll"

# Widen column B on Items to fit the new, longer item name.
$wsItems.Columns.Item(2).ColumnWidth = 30.65

# --- Selections: restore each sheet's own cursor position ---
$wsGlobalVars.Range("A22").Select()
$wsKeywords.Range("A10").Select()
$wsItems.Range("C8").Select()

# ObjectTexts becomes the active sheet/tab, selected last so it ends up the
# workbook's active tab (matches activeTab moving to ObjectTexts).
$wsObjectTexts.Range("C4").Select()
